# menu() is not working properly
# Rebuild the "items", "soldProduct" and "userAccount" sheets to match the
# corrected/expanded product & sales data, and restore the expected
# selection / active-sheet state.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    # Force a cell to be stored as text even when its content looks like a
    # number (keeps leading zeros such as phone numbers / PINs intact).
    param($ws, [string]$addr, [string]$val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# ---------------------------------------------------------------------
# Sheet "items"
# ---------------------------------------------------------------------
$items = $wb.Worksheets.Item("items")

# Row 2 & 3: product names/prices were swapped and corrected.
$items.Range("B2").Value = "Pen"
$items.Range("C2").Value = 10
$items.Range("D2").Value = 200

$items.Range("B3").Value = "Computer"
$items.Range("C3").Value = 10000
$items.Range("D3").Value = 10

# Rows 6-11: re-numbered / re-priced, plus one new row (Hard disk).
$items.Range("A6").Value = 1007
$items.Range("B6").Value = "Ipad"
$items.Range("C6").Value = 1000
$items.Range("D6").Value = 5

$items.Range("A7").Value = 1008
$items.Range("B7").Value = "mouse"
$items.Range("C7").Value = 300
$items.Range("D7").Value = 20

$items.Range("A8").Value = 1009
$items.Range("B8").Value = "Keyboard"
$items.Range("C8").Value = 100
$items.Range("D8").Value = 10

$items.Range("A9").Value = 1011
$items.Range("B9").Value = "Mouse"
$items.Range("C9").Value = 101
$items.Range("D9").Value = 8

$items.Range("A10").Value = 1012
$items.Range("B10").Value = "Smart watch"
$items.Range("C10").Value = 2000
$items.Range("D10").Value = 2

$items.Range("A11").Value = 1013
$items.Range("B11").Value = "Hard disk"
$items.Range("C11").Value = 8000
$items.Range("D11").Value = 5

# Rows 12-17: brand-new products appended at the bottom of the list.
$items.Range("A12").Value = 1014
$items.Range("B12").Value = "SSD"
$items.Range("C12").Value = 2800
$items.Range("D12").Value = 25

$items.Range("A13").Value = 1015
$items.Range("B13").Value = "RAM"
$items.Range("C13").Value = 3500
$items.Range("D13").Value = 25

$items.Range("A14").Value = 1016
$items.Range("B14").Value = "Monitor"
$items.Range("C14").Value = 15000
$items.Range("D14").Value = 5

$items.Range("A15").Value = 1017
$items.Range("B15").Value = "Laptop table"
$items.Range("C15").Value = 1000
$items.Range("D15").Value = 10

$items.Range("A16").Value = 1018
$items.Range("B16").Value = "Cable"
$items.Range("C16").Value = 100
$items.Range("D16").Value = 99950

$items.Range("A17").Value = 1019
$items.Range("B17").Value = "Wireless mouse"
$items.Range("C17").Value = 1000
$items.Range("D17").Value = 8

# ---------------------------------------------------------------------
# Sheet "soldProduct" - newly recorded sales appended after row 3.
# ---------------------------------------------------------------------
$sold = $wb.Worksheets.Item("soldProduct")

$sold.Range("A4").Value = 1011
$sold.Range("B4").Value = "Mouse"
$sold.Range("C4").Value = 8
$sold.Range("D4").Value = 101
$sold.Range("E4").Value = 808
$sold.Range("F4").Value = "Barkatopu"
$sold.Range("G4").Value = "abc"
Set-TextCell $sold "H4" "0180"
$sold.Range("I4").Value = "COD"

$sold.Range("A5").Value = 1012
$sold.Range("B5").Value = "Smart watch"
$sold.Range("C5").Value = 3
$sold.Range("D5").Value = 2000
$sold.Range("E5").Value = 6000
$sold.Range("F5").Value = "dip"
$sold.Range("G5").Value = "fftgr"
Set-TextCell $sold "H5" "34435"
$sold.Range("I5").Value = "COD"

$sold.Range("A6").Value = 1012
$sold.Range("B6").Value = "Smart watch"
$sold.Range("C6").Value = 3
$sold.Range("D6").Value = 2000
$sold.Range("E6").Value = 6000
$sold.Range("F6").Value = "Barkatopu"
$sold.Range("G6").Value = "abc"
Set-TextCell $sold "H6" "0180"
$sold.Range("I6").Value = "COD"

$sold.Range("A7").Value = 1013
Set-TextCell $sold "B7" "0"
$sold.Range("C7").Value = 1
$sold.Range("D7").Value = 0
$sold.Range("E7").Value = 0
$sold.Range("F7").Value = "Barkatopu"
$sold.Range("G7").Value = "abc"
Set-TextCell $sold "H7" "0180"
$sold.Range("I7").Value = "COD"

$sold.Range("A8").Value = 1018
$sold.Range("B8").Value = "Cable"
$sold.Range("C8").Value = 50
$sold.Range("D8").Value = 100
$sold.Range("E8").Value = 5000
$sold.Range("F8").Value = "Name"
$sold.Range("G8").Value = "Address"
$sold.Range("H8").Value = "Phone"
$sold.Range("I8").Value = "COD"

$sold.Range("A9").Value = 1019
$sold.Range("B9").Value = "Wireless mouse"
$sold.Range("C9").Value = 2
$sold.Range("D9").Value = 1000
$sold.Range("E9").Value = 2000
$sold.Range("F9").Value = "Name"
$sold.Range("G9").Value = "Address"
$sold.Range("H9").Value = "Phone"
$sold.Range("I9").Value = "COD"

# ---------------------------------------------------------------------
# Sheet "userAccount" - extra columns + the real user list.
# ---------------------------------------------------------------------
$user = $wb.Worksheets.Item("userAccount")

$user.Range("B1").Value = "User ID"
$user.Range("D1").Value = "Address"
$user.Range("E1").Value = "Phone"

$user.Range("A2").Value = "Mohammad Barkatullah "
$user.Range("B2").Value = "barkatopu"
Set-TextCell $user "C2" "1234"
$user.Range("D2").Value = "abc"
Set-TextCell $user "E2" "01521"

$user.Range("A3").Value = "Barkatopu"
$user.Range("B3").Value = "barkat1345"
Set-TextCell $user "C3" "1234"
$user.Range("D3").Value = "abc"
Set-TextCell $user "E3" "0180"

$user.Range("A4").Value = "Mohammad "
$user.Range("B4").Value = "barkat"
Set-TextCell $user "C4" "12345"
$user.Range("D4").Value = "foolan"
Set-TextCell $user "E4" "000000"

$user.Range("A5").Value = " Admin"
$user.Range("B5").Value = "Admin"
Set-TextCell $user "C5" "1234"
$user.Range("D5").Value = "xyz"
$user.Range("E5").Value = 0

$user.Range("A6").Value = "Barkat"
$user.Range("B6").Value = "boss"
Set-TextCell $user "C6" "1234"
$user.Range("D6").Value = "xzy"
$user.Range("E6").Value = 0

$user.Range("A7").Value = "partho"
$user.Range("B7").Value = "partho123"
Set-TextCell $user "C7" "123"
$user.Range("D7").Value = "fds"
$user.Range("E7").Value = 304585

$user.Range("A8").Value = "Alma"
$user.Range("B8").Value = "alma"
Set-TextCell $user "C8" "123"
$user.Range("D8").Value = "xas"
$user.Range("E8").Value = 2432

# ---------------------------------------------------------------------
# Restore the selection on each sheet, then activate "userAccount" last
# so it becomes the active/visible tab (matches the saved workbook view).
# ---------------------------------------------------------------------
$items.Range("B9").Select() | Out-Null
$sold.Range("J4").Select() | Out-Null
$user.Range("B8").Select() | Out-Null
$user.Activate() | Out-Null
